# Manish Pandey.xlsx — "complate!!-> scrapping whole ipl"
#
# 1. Rename the (only) sheet from "Sheet1" to "Manish Pandey".
# 2. Insert a new first column "matchNo" (shifting teamName..result one
#    column to the right, B..M instead of A..L).
# 3. Replace the single sample data row with the full eight-match scrape.
#
# All data cells in the source file are stored as TEXT (even the
# numeric-looking ones such as runs/balls/sr), so every write below forces
# the cell to Text (NumberFormat "@") before assigning the value, then
# clears the format again so no stray number-format style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Manish Pandey"

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# header row + 8 match rows, columns A..M (1-based row/col numbers already)
$data = @(
    @("matchNo", "teamName",            "batterName",    "states",                  "runs", "balls", "fours", "sixes", "sr",     "opponentTeamName",            "venue",       "date",          "result"),
    @("28th",    "Sunrisers Hyderabad", "Manish Pandey", "b Mustafizur Rahman",      "31",   "20",    "3",     "2",     "155.00", "Rajasthan Royals",            "Delhi",       "May 02",        "Royals won by 55 runs"),
    @("37th",    "Sunrisers Hyderabad", "Manish Pandey", "b Ravi Bishnoi",           "13",   "23",    "1",     "0",     "56.52",  "Punjab Kings",                "Sharjah",     "September 25",  "Punjab Kings won by 5 runs"),
    @("33rd",    "Sunrisers Hyderabad", "Manish Pandey", "c & b Rabada",             "17",   "16",    "1",     "0",     "106.25", "Delhi Capitals",              "Dubai (DSC)", "September 22",  "Capitals won by 8 wickets (with 13 balls remaining)"),
    @("23rd",    "Sunrisers Hyderabad", "Manish Pandey", "c du Plessis b Ngidi",     "61",   "46",    "5",     "1",     "132.60", "Chennai Super Kings",         "Delhi",       "April 28",      "Super Kings won by 7 wickets (with 9 balls remaining)"),
    @("3rd",     "Sunrisers Hyderabad", "Manish Pandey", "",                         "61",   "44",    "2",     "3",     "138.63", "Kolkata Knight Riders",       "Chennai",     "April 11",      "KKR won by 10 runs"),
    @("6th",     "Sunrisers Hyderabad", "Manish Pandey", "c Patel b Shahbaz Ahmed",  "38",   "39",    "2",     "2",     "97.43",  "Royal Challengers Bangalore", "Chennai",     "April 14",      "RCB won by 6 runs"),
    @("9th",     "Sunrisers Hyderabad", "Manish Pandey", "c Pollard b Chahar",       "2",    "7",     "0",     "0",     "28.57",  "Mumbai Indians",              "Chennai",     "April 17",      "Mumbai won by 13 runs")
)

$rowCount = $data.Length
for ($ri = 0; $ri -lt $rowCount; $ri++) {
    $rowValues = $data[$ri]
    $rowNum = $ri + 1
    $colCount = $rowValues.Length
    for ($ci = 0; $ci -lt $colCount; $ci++) {
        $colNum = $ci + 1
        $val = $rowValues[$ci]
        if ($val -ne "") {
            Set-TextValue $rowNum $colNum $val
        }
    }
}

# D6 ("states" for the 3rd match) is an explicit empty string in the
# source data (not simply a blank/untouched cell). A plain "" assignment
# clears/blanks a cell in Excel, so use the classic leading-apostrophe
# trick ('force text, no visible content') to get a real Text cell whose
# value is the empty string.
$d6 = $ws.Range("D6")
$d6.Value = "'"
$d6.ClearFormats()
